# My selenium Project CRM final commit
#
# - clientdata: "No record found." message in I2 is replaced by a new
#   "The record has been deleted." message, and the old text is pushed
#   down to I3 (column I is widened to fit the new, longer text).
# - projectdata: sample row values "abc" / "abc1" are replaced with more
#   meaningful "projectabc" / "projectabc1" placeholders.
# - projectdata becomes the active / selected sheet (was notedata before).

$wb = $excel.ActiveWorkbook

$wsNote    = $wb.Worksheets.Item("notedata")
$wsClient  = $wb.Worksheets.Item("clientdata")
$wsProject = $wb.Worksheets.Item("projectdata")

# ---- clientdata sheet ----
$wsClient.Range("I2").Value = "The record has been deleted."
$wsClient.Range("I3").Value = "No record found."
$wsClient.Columns.Item(9).ColumnWidth = 26.6
[void]$wsClient.Range("I2").Select()

# ---- notedata sheet: selection moves down to E2 ----
[void]$wsNote.Activate()
[void]$wsNote.Range("E2").Select()

# ---- projectdata sheet: update sample values and make it the active tab ----
$wsProject.Range("A2").Value = "projectabc"
$wsProject.Range("E2").Value = "projectabc1"
[void]$wsProject.Activate()
[void]$wsProject.Range("E2").Select()
